# Make a waypoint object type.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Map")

# Add new waypoint row data
$ws.Range("E2").Value = 5
$ws.Range("H2").Value = 5

# Move the selection to match the new editing position
$ws.Range("F14").Select()
